# Apply Harvard case classification update to the stats sheet.
# - Row 1: swap headers for the "average_doctor" / "average_doctor_old" columns (BP/BQ)
#   so that BP now holds the *_old header and BQ holds the current header.
# - Rows 4-13: update the underlying precision/recall/F-score/etc. statistics with
#   newly recomputed values, including the corresponding BP/BQ (average_doctor /
#   average_doctor_old) summary columns whose contents shift one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("AI4").Value = 0.229
$ws.Range("AJ4").Value = 0.067
$ws.Range("AU4").Value = 0.153
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.168
$ws.Range("BA4").Value = 1.908
$ws.Range("BB4").Value = 0.163
$ws.Range("BC4").Value = 0.403
$ws.Range("BG4").Value = 0.709
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.68
$ws.Range("BN4").Value = 0.082
$ws.Range("BO4").Value = 0.286
$ws.Range("BP4").Value = 0.636
$ws.Range("BQ4").Value = 0.644
$ws.Range("E4").Value = 0.408
$ws.Range("F4").Value = 0.07199999999999999
$ws.Range("G4").Value = 0.267
$ws.Range("N4").Value = 0.403
$ws.Range("O4").Value = 0.062
$ws.Range("P4").Value = 0.25
$ws.Range("W4").Value = 0.225
$ws.Range("X4").Value = 0.106
$ws.Range("Y4").Value = 0.325
$ws.Range("AI5").Value = 0.274
$ws.Range("AK5").Value = 0.313
$ws.Range("AU5").Value = 0.311
$ws.Range("AV5").Value = 0.104
$ws.Range("AW5").Value = 0.323
$ws.Range("BA5").Value = 1.369
$ws.Range("BB5").Value = 0.083
$ws.Range("BC5").Value = 0.287
$ws.Range("BG5").Value = 0.407
$ws.Range("BH5").Value = 0.056
$ws.Range("BI5").Value = 0.236
$ws.Range("BM5").Value = 0.586
$ws.Range("BN5").Value = 0.078
$ws.Range("BO5").Value = 0.279
$ws.Range("BP5").Value = 0.456
$ws.Range("BQ5").Value = 0.451
$ws.Range("E5").Value = 0.547
$ws.Range("F5").Value = 0.09
$ws.Range("G5").Value = 0.3
$ws.Range("N5").Value = 0.754
$ws.Range("O5").Value = 0.079
$ws.Range("P5").Value = 0.282
$ws.Range("W5").Value = 0.227
$ws.Range("X5").Value = 0.11
$ws.Range("Y5").Value = 0.332
$ws.Range("AI6").Value = 0.249
$ws.Range("AU6").Value = 0.205
$ws.Range("BA6").Value = 1.583
$ws.Range("BG6").Value = 0.517
$ws.Range("BM6").Value = 0.63
$ws.Range("BP6").Value = 0.528
$ws.Range("BQ6").Value = 0.527
$ws.Range("E6").Value = 0.467
$ws.Range("N6").Value = 0.525
$ws.Range("W6").Value = 0.226
$ws.Range("AI7").Value = 0.264
$ws.Range("AU7").Value = 0.258
$ws.Range("BA7").Value = 1.446
$ws.Range("BG7").Value = 0.445
$ws.Range("BM7").Value = 0.603
$ws.Range("BP7").Value = 0.482
$ws.Range("BQ7").Value = 0.479
$ws.Range("E7").Value = 0.512
$ws.Range("N7").Value = 0.642
$ws.Range("W7").Value = 0.227
$ws.Range("AI8").Value = 0.25
$ws.Range("AJ8").Value = 0.098
$ws.Range("AK8").Value = 0.313
$ws.Range("AU8").Value = 0.243
$ws.Range("AV8").Value = 0.073
$ws.Range("AW8").Value = 0.271
$ws.Range("BA8").Value = 1.681
$ws.Range("BG8").Value = 0.542
$ws.Range("BH8").Value = 0.104
$ws.Range("BI8").Value = 0.323
$ws.Range("BM8").Value = 0.697
$ws.Range("BN8").Value = 0.068
$ws.Range("BO8").Value = 0.261
$ws.Range("BP8").Value = 0.5600000000000001
$ws.Range("BQ8").Value = 0.574
$ws.Range("E8").Value = 0.576
$ws.Range("F8").Value = 0.113
$ws.Range("G8").Value = 0.336
$ws.Range("N8").Value = 0.761
$ws.Range("O8").Value = 0.068
$ws.Range("P8").Value = 0.26
$ws.Range("W8").Value = 0.235
$ws.Range("X8").Value = 0.114
$ws.Range("Y8").Value = 0.337
$ws.Range("AI9").Value = 0.14
$ws.Range("AJ9").Value = 0.12
$ws.Range("AK9").Value = 0.347
$ws.Range("BA9").Value = 1.582
$ws.Range("BB9").Value = 0.243
$ws.Range("BC9").Value = 0.493
$ws.Range("BG9").Value = 0.5580000000000001
$ws.Range("BH9").Value = 0.247
$ws.Range("BI9").Value = 0.497
$ws.Range("BM9").Value = 0.605
$ws.Range("BN9").Value = 0.239
$ws.Range("BO9").Value = 0.489
$ws.Range("BP9").Value = 0.527
$ws.Range("BQ9").Value = 0.534
$ws.Range("E9").Value = 0.488
$ws.Range("F9").Value = 0.25
$ws.Range("G9").Value = 0.5
$ws.Range("N9").Value = 0.651
$ws.Range("O9").Value = 0.227
$ws.Range("P9").Value = 0.477
$ws.Range("W9").Value = 0.14
$ws.Range("X9").Value = 0.12
$ws.Range("Y9").Value = 0.347
$ws.Range("AI10").Value = 0.279
$ws.Range("AJ10").Value = 0.201
$ws.Range("AK10").Value = 0.449
$ws.Range("AU10").Value = 0.233
$ws.Range("AV10").Value = 0.178
$ws.Range("AW10").Value = 0.422
$ws.Range("BA10").Value = 1.93
$ws.Range("BB10").Value = 0.249
$ws.Range("BC10").Value = 0.499
$ws.Range("BG10").Value = 0.605
$ws.Range("BH10").Value = 0.239
$ws.Range("BI10").Value = 0.489
$ws.Range("BM10").Value = 0.86
$ws.Range("BN10").Value = 0.12
$ws.Range("BO10").Value = 0.347
$ws.Range("BP10").Value = 0.643
$ws.Range("BQ10").Value = 0.679
$ws.Range("E10").Value = 0.628
$ws.Range("F10").Value = 0.234
$ws.Range("G10").Value = 0.483
$ws.Range("N10").Value = 0.86
$ws.Range("O10").Value = 0.12
$ws.Range("P10").Value = 0.347
$ws.Range("W10").Value = 0.279
$ws.Range("X10").Value = 0.201
$ws.Range("Y10").Value = 0.449
$ws.Range("AI11").Value = 0.279
$ws.Range("AJ11").Value = 0.201
$ws.Range("AK11").Value = 0.449
$ws.Range("AU11").Value = 0.326
$ws.Range("AV11").Value = 0.22
$ws.Range("AW11").Value = 0.469
$ws.Range("BA11").Value = 1.93
$ws.Range("BB11").Value = 0.249
$ws.Range("BC11").Value = 0.499
$ws.Range("BG11").Value = 0.605
$ws.Range("BH11").Value = 0.239
$ws.Range("BI11").Value = 0.489
$ws.Range("BM11").Value = 0.86
$ws.Range("BN11").Value = 0.12
$ws.Range("BO11").Value = 0.347
$ws.Range("BP11").Value = 0.643
$ws.Range("BQ11").Value = 0.679
$ws.Range("E11").Value = 0.651
$ws.Range("F11").Value = 0.227
$ws.Range("G11").Value = 0.477
$ws.Range("N11").Value = 0.884
$ws.Range("O11").Value = 0.103
$ws.Range("P11").Value = 0.321
$ws.Range("W11").Value = 0.279
$ws.Range("X11").Value = 0.201
$ws.Range("Y11").Value = 0.449
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.571
$ws.Range("AV12").Value = 1.673
$ws.Range("AW12").Value = 1.294
$ws.Range("BA12").Value = 3.632
$ws.Range("BB12").Value = 0.227
$ws.Range("BC12").Value = 0.477
$ws.Range("BG12").Value = 1.077
$ws.Range("BH12").Value = 0.07099999999999999
$ws.Range("BI12").Value = 0.266
$ws.Range("BM12").Value = 1.405
$ws.Range("BN12").Value = 0.457
$ws.Range("BO12").Value = 0.676
$ws.Range("BP12").Value = 1.211
$ws.Range("BQ12").Value = 1.283
$ws.Range("E12").Value = 1.429
$ws.Range("F12").Value = 0.673
$ws.Range("G12").Value = 0.821
$ws.Range("N12").Value = 1.513
$ws.Range("O12").Value = 1.07
$ws.Range("P12").Value = 1.035
$ws.Range("W12").Value = 1.75
$ws.Range("X12").Value = 0.6879999999999999
$ws.Range("Y12").Value = 0.829
$ws.Range("AI13").Value = 1.383
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK13").Value = 0.633
$ws.Range("AU13").Value = 2.475
$ws.Range("AV13").Value = 1.492
$ws.Range("AW13").Value = 1.222
$ws.Range("BA13").Value = 2.585
$ws.Range("BB13").Value = 0.316
$ws.Range("BC13").Value = 0.5629999999999999
$ws.Range("BG13").Value = 0.645
$ws.Range("BH13").Value = 0.096
$ws.Range("BI13").Value = 0.31
$ws.Range("BM13").Value = 1.022
$ws.Range("BN13").Value = 0.398
$ws.Range("BO13").Value = 0.631
$ws.Range("BP13").Value = 0.862
$ws.Range("BQ13").Value = 0.796
$ws.Range("E13").Value = 1.715
$ws.Range("F13").Value = 0.952
$ws.Range("G13").Value = 0.975
$ws.Range("N13").Value = 2.312
$ws.Range("O13").Value = 1.223
$ws.Range("P13").Value = 1.106
$ws.Range("W13").Value = 1.091
$ws.Range("X13").Value = 0.181
$ws.Range("Y13").Value = 0.426
